# Add a new page link row into the table on the active sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "/moderator/users"
$ws.Range("A10").Value = "ModeratorPanelUsersPage"
$ws.Range("C10").Value = "Anton Tsvihun"

# Match the selection state recorded after the edit
$ws.Range("C10").Select()
